$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.541.72"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "3.267.65"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.41"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "3.837.09"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.56"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "68.510.99"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").Value = "3.264.64"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "394.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.81%  "
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000119"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  +3.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "164.45"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.32%  "
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.36"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  -4.31%  "
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.28"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "346.74"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0688"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "2.611.69"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("E51").Value = "  -0.04%  "
